# [UPD] Hip12 + Facebook sort
#
# Applies the HIP-12 title-banner edit to the single slide in the deck:
#   - "Rectangle 4" banner text "HIP-12" -> "VOTE: HIP-12", shrinks the
#     font (130pt -> 96pt) and is repositioned/resized to the new
#     autofit box.
#   - "Rectangle 5" footer link text "talk.harmony.one" -> "gov.harmony.one".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Rectangle 4": the big "HIP-12" banner -> "VOTE: HIP-12" ---
$banner = $s.Shapes.Item(4)
$bannerRun = $banner.TextFrame.TextRange.Runs(1)
$bannerRun.Text = "VOTE: HIP-12"
$banner.TextFrame.TextRange.Font.Size = 96

# Reposition/resize the shape box to match the re-flowed text extent
# (values below are the EMU offsets/extents from the target layout,
# expressed in points as the Shape.Left/Top/Width/Height API expects).
$banner.Left = -2.4
$banner.Top = 7.197165354330709
$banner.Width = 957.75
$banner.Height = 123.59527559055118

# --- "Rectangle 5": the footer link "talk.harmony.one" -> "gov.harmony.one" ---
$link = $s.Shapes.Item(5)
$link.TextFrame.TextRange.Runs(1).Text = "gov.harmony.one"
